$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 17:25:46"
$wsZhCn.Range("H2").Value = "2016-03-24 17:26:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 17:25:51"
$wsDeDe.Range("H2").Value = "2016-03-24 17:26:37"
